$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A35").Value = "2025-04-28 23:50:10"
$ws.Range("B35").Value = 178
